$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(92, 8).Value = 47893464
$ws.Cells.Item(92, 9).Value = 2646173
$ws.Cells.Item(92, 11).Value = 2646173
$ws.Cells.Item(92, 13).Value = -2644925

$ws.Cells.Item(129, 8).Value = 877.0599999999999
$ws.Cells.Item(129, 9).Value = 465.25
$ws.Cells.Item(129, 10).Value = 980.0125
$ws.Cells.Item(129, 11).Value = 1395.75
$ws.Cells.Item(129, 12).Value = 2940.0375
$ws.Cells.Item(129, 13).Value = 3604.25
$ws.Cells.Item(129, 14).Value = -12940.0375

$ws.Cells.Item(132, 8).Value = 6411374
$ws.Cells.Item(132, 9).Value = 745.9211
$ws.Cells.Item(132, 10).Value = 23811650
$ws.Cells.Item(132, 11).Value = 2237.7633
$ws.Cells.Item(132, 12).Value = 71434950
$ws.Cells.Item(132, 13).Value = 292.2366999999999
$ws.Cells.Item(132, 14).Value = -71440010

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 4887.7866
$ws.Cells.Item(32, 9).Value = 3367.2532
$ws.Cells.Item(32, 11).Value = 3367.2532
$ws.Cells.Item(32, 13).Value = -3080.2532

$ws.Cells.Item(45, 8).Value = 6056.25
$ws.Cells.Item(45, 9).Value = 7275.9375
$ws.Cells.Item(45, 10).Value = 1177.5
$ws.Cells.Item(45, 11).Value = 7275.9375
$ws.Cells.Item(45, 12).Value = 1177.5
$ws.Cells.Item(45, 13).Value = -6898.9375
$ws.Cells.Item(45, 14).Value = -1931.5

$ws.Cells.Item(122, 8).Value = 916847.4399999999
$ws.Cells.Item(122, 9).Value = 1425425.4
$ws.Cells.Item(122, 10).Value = 1407.3
$ws.Cells.Item(122, 11).Value = 4276276.199999999
$ws.Cells.Item(122, 12).Value = 4221.9
$ws.Cells.Item(122, 13).Value = -4273826.199999999
$ws.Cells.Item(122, 14).Value = -9121.9

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(35, 8).Value = 23400
$ws.Cells.Item(35, 10).Value = 23400
$ws.Cells.Item(35, 12).Value = 23400
$ws.Cells.Item(35, 14).Value = -24020

$ws.Cells.Item(94, 8).Value = 2130.8635
$ws.Cells.Item(94, 10).Value = 2597.7778
$ws.Cells.Item(94, 12).Value = 2597.7778
$ws.Cells.Item(94, 14).Value = -3499.7778

$ws.Cells.Item(134, 8).Value = 3706.463
$ws.Cells.Item(134, 9).Value = 4355.027
$ws.Cells.Item(134, 10).Value = 2294.8823
$ws.Cells.Item(134, 11).Value = 13065.081
$ws.Cells.Item(134, 12).Value = 6884.646900000001
$ws.Cells.Item(134, 13).Value = -10530.081
$ws.Cells.Item(134, 14).Value = -11954.6469

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 1531.125
$ws.Cells.Item(16, 9).Value = 1064
$ws.Cells.Item(16, 10).Value = 1894.4445
$ws.Cells.Item(16, 11).Value = 1064
$ws.Cells.Item(16, 12).Value = 1894.4445
$ws.Cells.Item(16, 13).Value = -777
$ws.Cells.Item(16, 14).Value = -2468.4445

$ws.Cells.Item(31, 8).Value = 11366384
$ws.Cells.Item(31, 9).Value = 1032.0344
$ws.Cells.Item(31, 10).Value = 33339398
$ws.Cells.Item(31, 11).Value = 1032.0344
$ws.Cells.Item(31, 12).Value = 33339398
$ws.Cells.Item(31, 13).Value = -737.0344
$ws.Cells.Item(31, 14).Value = -33339988

$ws.Cells.Item(34, 8).Value = 11366384
$ws.Cells.Item(34, 9).Value = 1032.0344
$ws.Cells.Item(34, 10).Value = 33339398
$ws.Cells.Item(34, 11).Value = 1032.0344
$ws.Cells.Item(34, 12).Value = 33339398
$ws.Cells.Item(34, 13).Value = -830.0344
$ws.Cells.Item(34, 14).Value = -33339802

$ws.Cells.Item(58, 8).Value = 6668096
$ws.Cells.Item(58, 9).Value = 9260270
$ws.Cells.Item(58, 10).Value = 2506
$ws.Cells.Item(58, 11).Value = 9260270
$ws.Cells.Item(58, 12).Value = 2506
$ws.Cells.Item(58, 13).Value = -9260067
$ws.Cells.Item(58, 14).Value = -2912

$ws.Cells.Item(94, 8).Value = 4025.2727
$ws.Cells.Item(94, 9).Value = 4439.5557
$ws.Cells.Item(94, 10).Value = 3738.4614
$ws.Cells.Item(94, 11).Value = 4439.5557
$ws.Cells.Item(94, 12).Value = 3738.4614
$ws.Cells.Item(94, 13).Value = -3988.5557
$ws.Cells.Item(94, 14).Value = -4640.4614

$ws.Cells.Item(113, 8).Value = 1531.125
$ws.Cells.Item(113, 9).Value = 1064
$ws.Cells.Item(113, 10).Value = 1894.4445
$ws.Cells.Item(113, 11).Value = 1064
$ws.Cells.Item(113, 12).Value = 1894.4445
$ws.Cells.Item(113, 13).Value = 1106
$ws.Cells.Item(113, 14).Value = -6234.4445

$ws.Cells.Item(134, 8).Value = 8131459.5
$ws.Cells.Item(134, 9).Value = 13334845
$ws.Cells.Item(134, 10).Value = 1169
$ws.Cells.Item(134, 11).Value = 40004535
$ws.Cells.Item(134, 12).Value = 3507
$ws.Cells.Item(134, 13).Value = -40002000
$ws.Cells.Item(134, 14).Value = -8577

$ws.Cells.Item(136, 8).Value = 6668096
$ws.Cells.Item(136, 9).Value = 9260270
$ws.Cells.Item(136, 10).Value = 2506
$ws.Cells.Item(136, 11).Value = 27780810
$ws.Cells.Item(136, 12).Value = 7518
$ws.Cells.Item(136, 13).Value = -27778260
$ws.Cells.Item(136, 14).Value = -12618

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(33, 8).Value = 4550065
$ws.Cells.Item(33, 9).Value = 5882417.5
$ws.Cells.Item(33, 10).Value = 20066
$ws.Cells.Item(33, 11).Value = 35294505
$ws.Cells.Item(33, 12).Value = 120396
$ws.Cells.Item(33, 13).Value = -35294222
$ws.Cells.Item(33, 14).Value = -120962

$ws.Cells.Item(75, 8).Value = 23809990
$ws.Cells.Item(75, 9).Value = 199.75
$ws.Cells.Item(75, 10).Value = 71429570
$ws.Cells.Item(75, 11).Value = 599.25
$ws.Cells.Item(75, 12).Value = 214288710
$ws.Cells.Item(75, 13).Value = 398.75
$ws.Cells.Item(75, 14).Value = -214290706

$ws.Cells.Item(78, 8).Value = 23809990
$ws.Cells.Item(78, 9).Value = 199.75
$ws.Cells.Item(78, 10).Value = 71429570
$ws.Cells.Item(78, 11).Value = 1797.75
$ws.Cells.Item(78, 12).Value = 642866130
$ws.Cells.Item(78, 13).Value = 3194.25
$ws.Cells.Item(78, 14).Value = -642876114

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(126, 8).Value = 9223.214
$ws.Cells.Item(126, 9).Value = 12052.4
$ws.Cells.Item(126, 10).Value = 2150.25
$ws.Cells.Item(126, 11).Value = 36157.2
$ws.Cells.Item(126, 12).Value = 6450.75
$ws.Cells.Item(126, 13).Value = -33687.2
$ws.Cells.Item(126, 14).Value = -11390.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 4446234.5
$ws.Cells.Item(22, 9).Value = 12346429
$ws.Cells.Item(22, 10).Value = 2374.5
$ws.Cells.Item(22, 11).Value = 12346429
$ws.Cells.Item(22, 12).Value = 2374.5
$ws.Cells.Item(22, 13).Value = -12346134
$ws.Cells.Item(22, 14).Value = -2964.5

$ws.Cells.Item(27, 8).Value = 4446234.5
$ws.Cells.Item(27, 9).Value = 12346429
$ws.Cells.Item(27, 10).Value = 2374.5
$ws.Cells.Item(27, 11).Value = 12346429
$ws.Cells.Item(27, 12).Value = 2374.5
$ws.Cells.Item(27, 13).Value = -12346322
$ws.Cells.Item(27, 14).Value = -2588.5

$ws.Cells.Item(46, 8).Value = 10417294
$ws.Cells.Item(46, 9).Value = 55556000
$ws.Cells.Item(46, 11).Value = 55556000
$ws.Cells.Item(46, 13).Value = -55555812

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(81, 8).Value = 15387391
$ws.Cells.Item(81, 9).Value = 857
$ws.Cells.Item(81, 10).Value = 25003974
$ws.Cells.Item(81, 11).Value = 1714
$ws.Cells.Item(81, 12).Value = 50007948
$ws.Cells.Item(81, 13).Value = -653
$ws.Cells.Item(81, 14).Value = -50010070

$ws.Cells.Item(84, 8).Value = 15387391
$ws.Cells.Item(84, 9).Value = 857
$ws.Cells.Item(84, 10).Value = 25003974
$ws.Cells.Item(84, 11).Value = 8570
$ws.Cells.Item(84, 12).Value = 250039740
$ws.Cells.Item(84, 13).Value = -3266
$ws.Cells.Item(84, 14).Value = -250050348

$ws.Cells.Item(132, 8).Value = 1115.1091
$ws.Cells.Item(132, 9).Value = 740.6512
$ws.Cells.Item(132, 10).Value = 2456.9167
$ws.Cells.Item(132, 11).Value = 2221.9536
$ws.Cells.Item(132, 12).Value = 7370.750100000001
$ws.Cells.Item(132, 13).Value = 308.0464000000002
$ws.Cells.Item(132, 14).Value = -12430.7501
